$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

function Set-TextValue {
    param($row, $col, $val)
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $val
}

Set-TextValue 2 4 "290.43"
Set-TextValue 2 5 "0.64%"
Set-TextValue 3 4 "30.76"
Set-TextValue 3 5 "-2.08%"
Set-TextValue 4 4 "4.875"
Set-TextValue 4 5 "-1.26%"
Set-TextValue 5 4 "0.07242"
Set-TextValue 5 5 "-0.58%"
Set-TextValue 6 4 "2.349"
Set-TextValue 6 5 "31.79%"
Set-TextValue 7 4 "7.650"
Set-TextValue 7 5 "-0.44%"
Set-TextValue 8 4 "3.704"
Set-TextValue 8 5 "-1.28%"
Set-TextValue 9 5 "-1.48%"
Set-TextValue 10 4 "0.1673"
Set-TextValue 10 5 "0.78%"
Set-TextValue 11 4 "0.08056"
Set-TextValue 11 5 "5.30%"
Set-TextValue 12 4 "0.08152"
Set-TextValue 12 5 "-0.91%"
Set-TextValue 13 4 "0.03073"
Set-TextValue 13 5 "1.66%"
Set-TextValue 14 4 "0.1002"
Set-TextValue 14 5 "-0.25%"
Set-TextValue 15 5 "-0.94%"
Set-TextValue 16 4 "0.005774"
Set-TextValue 16 5 "-1.54%"
Set-TextValue 17 4 "3.477"
Set-TextValue 17 5 "0.60%"
Set-TextValue 18 4 "2.077"
Set-TextValue 18 5 "-2.48%"
Set-TextValue 19 5 "1.59%"
Set-TextValue 20 5 "-0.65%"
Set-TextValue 21 4 "3.973"
Set-TextValue 21 5 "-9.35%"
Set-TextValue 22 4 "0.2109"
Set-TextValue 22 5 "5.45%"
Set-TextValue 23 4 "0.04534"
Set-TextValue 23 5 "0.85%"
Set-TextValue 24 4 "0.001214"
Set-TextValue 24 5 "-2.26%"
Set-TextValue 25 4 "0.004412"
Set-TextValue 25 5 "10.40%"
Set-TextValue 26 4 "0.0001301"
Set-TextValue 26 5 "2.63%"
Set-TextValue 27 4 "0.0003397"
Set-TextValue 27 5 "-95.47%"
Set-TextValue 39 5 "-1.91%"
Set-TextValue 40 4 "0.04378"
Set-TextValue 40 5 "-0.11%"
Set-TextValue 41 4 "0.007295"
Set-TextValue 41 5 "-2.99%"
Set-TextValue 43 5 "0.18%"
Set-TextValue 44 4 "0.002081"
Set-TextValue 44 5 "-10.76%"
Set-TextValue 45 4 "0.009176"
Set-TextValue 45 5 "-16.28%"
Set-TextValue 46 4 "0.00005703"
Set-TextValue 46 5 "-7.14%"
Set-TextValue 47 4 "0.00000000750"
Set-TextValue 47 5 "-1.19%"
Set-TextValue 48 4 "2.242"
Set-TextValue 48 5 "21.51%"
Set-TextValue 49 4 "0.002901"
Set-TextValue 49 5 "-4.49%"
Set-TextValue 50 4 "0.00002101"
Set-TextValue 50 5 "-1.19%"
Set-TextValue 51 4 "0.0002001"
Set-TextValue 51 5 "-1.19%"

Write-Output "Applied 70 price/volume updates to Sheet1."
